# Apply cryptos list price/volume updates (and a few row re-orderings)
# mirroring the upstream GitHub Actions scrape refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.365.56"
$ws.Range("E2").Value = "  +3.33%  "
$ws.Range("D3").Value = "3.787.66"
$ws.Range("E3").Value = "  +6.99%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.997"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "420.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.69%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.00"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.93%  "
$ws.Range("D7").Value = "3.950.42"
$ws.Range("E7").Value = "  +11.86%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.650"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.26%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.999"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.06%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.776"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.49%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.189"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +11.91%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000416"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +51.74%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "43.08"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.76%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "10.55"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.00%  "
$ws.Range("D15").Value = "4.361.52"
$ws.Range("E15").Value = "  +7.15%  "
$ws.Range("E16").Value = "  -0.56%  "
$ws.Range("D17").Value = "3.758.63"
$ws.Range("E17").Value = "  +6.75%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.52"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.21%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.17"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.77%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.13"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.97%  "
$ws.Range("D21").Value = "67.926.08"
$ws.Range("E21").Value = "  +3.15%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "446.27"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.88%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "15.62"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +18.33%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "90.41"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.09"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.26%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "38.46"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +11.76%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.40"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.09"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.35%  "
$ws.Range("E29").Value = "  +5.03%  "
$ws.Range("E30").Value = "  +4.96%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "12.63"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.86%  "
$ws.Range("E32").Value = "  +0.93%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.17"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.28%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.163"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "41.27"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.26%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "58.21"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.24%  "
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0489"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.24%  "
$ws.Range("B39").Value = "EnergySwap"
$ws.Range("C39").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "30.43"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +42.14%  "
$ws.Range("B40").Value = "PEPE"
$ws.Range("C40").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D40").Value = "0.0₃0712"
$ws.Range("E40").Value = "  -3.41%  "
$ws.Range("B41").Value = "Stellar"
$ws.Range("C41").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.148"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("B42").Value = "ThetaToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.98"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +28.60%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.993"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.43%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.38"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.21%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "147.83"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.24%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.18"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +24.03%  "
$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.91"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.70%  "
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.09"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.39%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.32"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.60"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -7.07%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.304"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.76%  "
